$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

Set-TextValue "B2" "0.25"
Set-TextValue "C2" "35041"
Set-TextValue "E2" "0.573"
Set-TextValue "F2" "15.788"
Set-TextValue "I2" "1702"
Set-TextValue "P2" "87.941"
Set-TextValue "Q2" "49759.031"
